# LOB1008.xlsx edit — reflects the commit's sharedStrings/sheet1 changes.
# Net effect (verified against the OOXML diff): several long Portuguese
# paragraphs were dropped from the shared-string table, the "Docentes
# responsaveis" / "Programa resumido" / "Programa" / "Avaliacao" /
# "Bibliografia" rows were re-populated with different (shorter/shifted)
# values, and the final row (old row 22, Bibliografia's long text) was
# removed outright, shifting the sheet's used range from C22 to C21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("Objetivos:") — body text replaced by the professor line.
$ws.Range("B10").Value = "6376612 - Daisy Rafaela da Silva"
$ws.Range("C10").Value = "6376612 - Daisy Rafaela da Silva"

# Row 13 — gains a label in A, B/C become "Semestral".
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14 — label becomes "Short syllabus:", B/C cleared entirely.
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14:C14").Clear()
$ws.Rows.Item(14).RowHeight = 60

# Row 15 — label becomes "Programa:", gains B/C = date string, height 60->120.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"
$ws.Rows.Item(15).RowHeight = 120

# Row 16 — label becomes "Syllabus:", B/C cleared entirely (height stays 120).
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16:C16").Clear()
$ws.Rows.Item(16).RowHeight = 120

# Row 17 — label becomes "Avaliacao:", loses its custom height (back to default).
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

# Row 18 — label becomes "Metodo:", gains B/C = professor line, height none->60.
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "6376612 - Daisy Rafaela da Silva"
$ws.Range("C18").Value = "6376612 - Daisy Rafaela da Silva"
$ws.Rows.Item(18).RowHeight = 60

# Row 19 — label becomes "Criterio:" (B/C text unchanged).
$ws.Range("A19").Value = "Critério:"

# Row 20 — label becomes "Norma de recuperacao:" (B/C text unchanged).
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21 — label becomes "Bibliografia:", height 60->120 (B/C text unchanged).
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# Row 22 (old Bibliografia long text) is dropped entirely; the sheet's used
# range shrinks from A1:C22 to A1:C21.
$ws.Rows.Item(22).Delete()

Write-Output "edit complete"
